$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update TERMINATION (col H) and NOTES (col I) labels to reflect the
# relabeled Arduino LCD backpack pinout.
$ws.Range("H14").Value = "ARD-A5"
$ws.Range("H15").Value = "ARD-A4"
$ws.Range("H16").Value = "ARD-D11"
$ws.Range("H17").Value = "ARD-D10"
$ws.Range("I22").Value = "Ard D7"
$ws.Range("I23").Value = "Ard D6"

# Update the last active selection shown when the sheet is reopened.
$ws.Range("P29").Select()
